$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the numbers stored as text (matches the source data's original
# "number stored as text" formatting) instead of letting Excel coerce
# them to numeric values.
$ws.Range("C2:F16").NumberFormat = "@"

# Target values for columns C (runs), D (balls), E (fours), F (sixes)
# for rows 2-16, reflecting the re-sorted / updated match data.
$data = @{
    2  = @(19, 20, 1, 1)
    3  = @(51, 38, 6, 2)
    4  = @(12, 11, 1, 0)
    5  = @(0, 4, 0, 0)
    6  = @(0, 2, 0, 0)
    7  = @(36, 29, 5, 0)
    8  = @(53, 32, 6, 1)
    9  = @(47, 28, 6, 1)
    10 = @(40, 26, 4, 1)
    11 = @(10, 7, 2, 0)
    12 = @(17, 16, 2, 0)
    13 = @(10, 10, 1, 0)
    14 = @(79, 47, 11, 2)
    15 = @(79, 43, 10, 3)
    16 = @(27, 18, 6, 0)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 3).Value = [string]$vals[0]
    $ws.Cells.Item($row, 4).Value = [string]$vals[1]
    $ws.Cells.Item($row, 5).Value = [string]$vals[2]
    $ws.Cells.Item($row, 6).Value = [string]$vals[3]
}
